# Front End Web Developer.xlsx - apply commit:
# "New MDN resources on Flexbox added. Some google and frontmaster courses added."
#
# Changes:
#  1. Insert a new blank "Sheet2" tab right after "Sheet1".
#  2. Insert a new "Frontend Master" tab right after "MDN".
#  3. Insert a new "CSS Tricks" tab right after "Frontend Master".
#  4. Add 3 new Flexbox/background-image rows to "MDN".
#  5. Add a new "Google Fonts" row to "Google".
#  6. Re-activate "Google" tab (it was the active tab before the edit too).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "Sheet2" tab, placed after "Sheet1"
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$sheet2.Name = "Sheet2"

$sheet2.Columns.Item(2).ColumnWidth = 61.666666666666664   # -> stored width 62.44140625 (nearest achievable)
$sheet2.Columns.Item(3).ColumnWidth = 70.33333333333333    # -> stored width 71.109375   (nearest achievable)

$sheet2.Range("B2").Value = "MIT Font"
$sheet2.Range("C2").Value = "https://web.mit.edu/jmorzins/www/fonts.html"
$sheet2.Range("C3").Select()

# ---------------------------------------------------------------------------
# 2 & 3. New "Frontend Master" and "CSS Tricks" tabs, placed after "MDN"
# ---------------------------------------------------------------------------
$mdnSheet = $wb.Worksheets.Item("MDN")
$frontendMaster = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $mdnSheet)
$frontendMaster.Name = "Frontend Master"

$cssTricks = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $frontendMaster)
$cssTricks.Name = "CSS Tricks"

# --- Frontend Master content ---
$frontendMaster.Columns.Item(2).ColumnWidth = 61.166666666666664  # -> stored width 62   (exact)
$frontendMaster.Columns.Item(3).ColumnWidth = 114.66666666666667  # -> stored width 115.5 (nearest achievable)

$fmB2 = $frontendMaster.Range("B2")
$fmB2.Value = "Firebase with React, v2 "
$fmB2.Font.Name = "Arial"
$fmB2.VerticalAlignment = -4108
$fmB2.WrapText = $true
$frontendMaster.Range("C2").Value = "https://frontendmasters.com/courses/firebase-react-v2"
$frontendMaster.Range("C2").Select()

# --- CSS Tricks content ---
$cssTricks.Columns.Item(2).ColumnWidth = 61.833333333333336  # -> stored width 62.6640625 (nearest achievable)
$cssTricks.Columns.Item(3).ColumnWidth = 79.16666666666667   # -> stored width 80         (exact)

$ctB2 = $cssTricks.Range("B2")
$ctB2.Value = "Firebase Crash Course "
$ctB2.Font.Name = "Lato"
$ctB2.Font.Bold = $true
$ctB2.VerticalAlignment = -4108
$ctB2.WrapText = $true
$cssTricks.Range("C2").Value = "https://css-tricks.com/firebase-crash-course/"
$cssTricks.Range("B4").Select()

# ---------------------------------------------------------------------------
# 4. New MDN rows (Flexbox resources)
# ---------------------------------------------------------------------------
$mdnSheet.Range("B12").Value = "Background Images"
$mdnSheet.Range("C12").Value = "https://developer.mozilla.org/en-US/docs/Web/CSS/background-image"
$mdnSheet.Range("B14").Value = "Flexbox"
$mdnSheet.Range("C14").Value = "https://developer.mozilla.org/en-US/docs/Learn/CSS/CSS_layout/Flexbox"
$mdnSheet.Range("B16").Value = "Flexbox Basic Concept"
$mdnSheet.Range("C16").Value = "https://developer.mozilla.org/en-US/docs/Web/CSS/CSS_Flexible_Box_Layout/Basic_Concepts_of_Flexbox"
$mdnSheet.Range("C18").Select()

# ---------------------------------------------------------------------------
# 5. New Google row (Google Fonts)
# ---------------------------------------------------------------------------
$googleSheet = $wb.Worksheets.Item("Google")
$googleSheet.Range("B6").Value = "Google Fonts"
$googleSheet.Range("C6").Value = "https://fonts.google.com/"

# ---------------------------------------------------------------------------
# 6. Keep "Google" as the active tab (matches pre-edit state)
# ---------------------------------------------------------------------------
$googleSheet.Range("B6").Select()
$googleSheet.Activate()
